$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (Price, Volume(1h)) updates. $null means the value is unchanged.
$updates = @{
    2  = @('27.518.42', '  +1.89%  ')
    3  = @('1.563.00', $null)
    4  = @('0.990', '  -1.53%  ')
    5  = @('210.60', '  +1.21%  ')
    6  = @('0.491', '  +0.20%  ')
    7  = @($null, '  -1.58%  ')
    8  = @('22.43', '  +1.65%  ')
    9  = @('0.251', '  +0.68%  ')
    10 = @($null, '  -0.08%  ')
    11 = @('0.0868', '  +1.44%  ')
    12 = @('1.786.32', '  +0.17%  ')
    13 = @('1.561.78', '  +0.18%  ')
    14 = @($null, '  +0.65%  ')
    15 = @($null, '  +0.38%  ')
    16 = @('27.502.92', '  +1.79%  ')
    17 = @('62.49', '  +1.11%  ')
    18 = @('224.31', '  +4.08%  ')
    19 = @('7.52', '  +1.71%  ')
    20 = @($null, '  +0.18%  ')
    22 = @($null, '  +0.23%  ')
    23 = @('9.40', '  +2.00%  ')
    25 = @('150.11', '  -2.01%  ')
    26 = @($null, '  +2.60%  ')
    27 = @($null, '  +0.40%  ')
    28 = @('15.16', '  +0.78%  ')
    29 = @($null, '  -1.22%  ')
    30 = @('1.14', '  +1.26%  ')
    31 = @($null, '  -0.69%  ')
    32 = @($null, '  +0.28%  ')
    33 = @('1.463.51', '  +2.93%  ')
    34 = @('3.18', '  -0.20%  ')
    35 = @('1.10', '  +2.51%  ')
    36 = @($null, '  +0.89%  ')
    37 = @('2.30', '  -1.63%  ')
    38 = @($null, '  +0.35%  ')
    39 = @($null, '  +2.08%  ')
    40 = @($null, '  +0.63%  ')
    41 = @('5.70', '  -1.67%  ')
    42 = @($null, '  +1.03%  ')
    43 = @('0.990', '  -1.51%  ')
    44 = @($null, '  +8.52%  ')
    45 = @('0.976', '  -2.15%  ')
    46 = @('65.00', '  +0.55%  ')
    47 = @('1.701.43', '  +0.26%  ')
    48 = @('86.71', '  +0.02%  ')
    49 = @($null, '  +0.41%  ')
    50 = @($null, '  +1.04%  ')
    51 = @('0.0952', '  -0.74%  ')
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $priceVal = $vals[0]
    $volVal = $vals[1]

    if ($null -ne $priceVal) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $priceVal
        $cell.Style = "Normal"
    }
    if ($null -ne $volVal) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $volVal
        $cell.Style = "Normal"
    }
}
